$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text so values like "583.67" are not
# auto-converted to numbers by Excel, matching the original inlineStr cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.813.83'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.252.28'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '583.67'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = '183.24'
$ws.Range('E6').Value = '  +3.65%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E9').Value = '  +4.09%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').Value = '3.817.01'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '28.56'
$ws.Range('E14').Value = '  +2.54%  '
$ws.Range('D15').Value = '67.826.80'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('D17').Value = '3.249.07'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '5.83'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = '13.57'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').Value = '381.12'
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('D21').Value = '7.64'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '71.23'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').Value = '0.513'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('D26').Value = '9.81'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '5.66'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = '22.88'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  +5.73%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +3.00%  '
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('D36').Value = '161.28'
$ws.Range('E36').Value = '  -6.17%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = '1.85'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = '0.835'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = '26.57'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '4.60'
$ws.Range('E40').Value = '  +7.23%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '6.68'
$ws.Range('E41').Value = '  +4.05%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('D43').Value = '41.27'
$ws.Range('E43').Value = '  +2.12%  '
$ws.Range('D44').Value = '347.05'
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').Value = '25.43'
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('D46').Value = '0.0687'
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('D47').Value = '2.633.67'
$ws.Range('E47').Value = '  -3.14%  '
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('E51').Value = '  +3.04%  '

# Restore the default (unstyled) appearance for the Price column so no stray
# cell style is introduced by the temporary text NumberFormat above.
$ws.Range("D2:D51").Style = "Normal"

